$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Pass/Fail" results column with "Pass" for every acceptance
# criteria row (C2:C7), matching the screenshot referenced in the commit.
$ws.Range("C2:C7").Value = "Pass"

# C2 already carried the plain bordered style (s=4); copy that same
# formatting down onto C3:C7 so the whole result column is consistent
# (previously C3:C7 used the wrap-text style inherited from columns A/B).
$ws.Range("C2").Copy()
$ws.Range("C3:C7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# Leave the selection where the user last left it.
$ws.Range("A9").Select()
